$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The three pasted screenshots in the document were inserted while the
#    proofing tools were running over them; mark each inline picture's run
#    as NoProofing so Word stamps <w:rPr><w:noProof/></w:rPr> on it (matches
#    the three +noProof hunks in the diff).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $d.InlineShapes.Item($i).Range.NoProofing = $true
}

# ---------------------------------------------------------------------------
# 2) Exercise C used a magic number (6) as the loop bound instead of the
#    num_cities variable - correct it: "i<6;" -> "i<num_cities;"
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*i<6;*") {
        $targetPara = $p
    }
}

$startP = $targetPara.Range.Start
$endP = $targetPara.Range.End

# Find the lone "6" inside that paragraph - use a freshly-created Range (not
# the paragraph's own live Range) so Find is properly scoped to it.
$pr = $d.Range($startP, $endP)
$pr.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

$sixStart = $pr.Start
$sixEnd = $pr.End

# Insert the replacement text right after the "6" run while its formatting
# still differs from the surrounding "variable name" colour - this keeps the
# new text in its own run instead of being silently merged back into the old
# run once the colours match.
$pr.InsertAfter("num_cities")

# Turn the original "6" character into a single space, coloured the same as
# the other variable names (E06C75), then colour the new "num_cities" text
# the same way - done as two separate operations (after the insertion above)
# so the two pieces stay in their own runs.
$rSpace = $d.Range($sixStart, $sixEnd)
$rSpace.Text = " "
$rSpace.Font.Color = 7695584

$rWord = $d.Range($sixEnd, $sixEnd + 10)
$rWord.Font.Color = 7695584
